$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at position 555; this pushes the existing
# rows 555:581 down to 556:582 and grows the used range to A1:R582.
$ws.Rows(555).Insert()

# Populate the newly inserted row 555 with the new weekly record.
$ws.Cells.Item(555, 1).Value = 4
$ws.Cells.Item(555, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(555, 3).Value = "Los Lagos"
$ws.Cells.Item(555, 4).Value = 44939
$ws.Cells.Item(555, 5).Value = 10
$ws.Cells.Item(555, 6).Value = 100114001
$ws.Cells.Item(555, 7).Value = "Papa"
$ws.Cells.Item(555, 8).Value = "Patagonia"
$ws.Cells.Item(555, 9).Value = "1a nueva(o)"
$ws.Cells.Item(555, 10).Value = 600
$ws.Cells.Item(555, 11).Value = 13000
$ws.Cells.Item(555, 12).Value = 14000
$ws.Cells.Item(555, 13).Value = 13500
$ws.Cells.Item(555, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(555, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(555, 16).Value = 540
$ws.Cells.Item(555, 17).Value = 25
$ws.Cells.Item(555, 18).Value = "Hortaliza"
